$d = $word.ActiveDocument

# Step 1: Replace the old note text with the new note text (single Find/Replace
# operation keeps the surrounding paragraph/bookmark intact, and the tracked
# "_GoBack" bookmark automatically follows to the end of the replaced text).
$oldText = "Note: the serial monitor cannot be used to read the string at the same time the 7-segment display is used to output data, you should use: Serial.begin(), then read the input, then use Serial.end() before displaying to the 7-segment display."
$newText = "Note: you can use Serial.available() to check whether a string is ready to be read in via the user."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Step 2: Locate the paragraph that now holds the new text.
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Note: you can use Serial.available()*") {
        $para = $candidate
        break
    }
}

$paraRange = $para.Range

# Step 3: Move the "_GoBack" bookmark so that it sits right before the final
# "." (i.e. between "...via the user" and "."), matching the target layout.
# Re-adding a bookmark at a new position naturally splits the run there.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$periodPos = $paraRange.End - 2
$bookmarkRange = $d.Range($periodPos, $periodPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Step 4: Split off the leading "Note: " into its own run by briefly adding
# and removing a zero-length bookmark at that boundary.
$noteSplitPos = $paraRange.Start + 6
$splitRange = $d.Range($noteSplitPos, $noteSplitPos)
$d.Bookmarks.Add("zzzTempSplit", $splitRange)
$d.Bookmarks.Item("zzzTempSplit").Delete()
